# Update MVP worksheet: rename Year->Season with full season labels,
# append basketball-reference player IDs to player names, and drop the
# now-redundant Award column (every row was "MVP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "Season"

# --- Season (col A) + Player (col B) data rows ---
$data = @(
    @("2017-18", "James Harden\hardeja01"),
    @("2016-17", "Russell Westbrook\westbru01"),
    @("2015-16", "Stephen Curry\curryst01"),
    @("2014-15", "Stephen Curry\curryst01"),
    @("2013-14", "Kevin Durant\duranke01"),
    @("2012-13", "LeBron James\jamesle01"),
    @("2011-12", "LeBron James\jamesle01"),
    @("2010-11", "Derrick Rose\rosede01"),
    @("2009-10", "LeBron James\jamesle01"),
    @("2008-09", "LeBron James\jamesle01"),
    @("2007-08", "Kobe Bryant\bryanko01"),
    @("2006-07", "Dirk Nowitzki\nowitdi01"),
    @("2005-06", "Steve Nash\nashst01"),
    @("2004-05", "Steve Nash\nashst01"),
    @("2003-04", "Kevin Garnett\garneke01"),
    @("2002-03", "Tim Duncan\duncati01"),
    @("2001-02", "Tim Duncan\duncati01"),
    @("2000-01", "Allen Iverson\iversal01"),
    @("1999-00", "Shaquille O'Neal\onealsh01"),
    @("1998-99", "Karl Malone\malonka01"),
    @("1997-98", "Michael Jordan\jordami01"),
    @("1996-97", "Karl Malone\malonka01"),
    @("1995-96", "Michael Jordan\jordami01"),
    @("1994-95", "David Robinson\robinda01"),
    @("1993-94", "Hakeem Olajuwon\olajuha01"),
    @("1992-93", "Charles Barkley\barklch01"),
    @("1991-92", "Michael Jordan\jordami01"),
    @("1990-91", "Michael Jordan\jordami01"),
    @("1989-90", "Magic Johnson\johnsma02"),
    @("1988-89", "Magic Johnson\johnsma02"),
    @("1987-88", "Michael Jordan\jordami01"),
    @("1986-87", "Magic Johnson\johnsma02"),
    @("1985-86", "Larry Bird\birdla01"),
    @("1984-85", "Larry Bird\birdla01"),
    @("1983-84", "Larry Bird\birdla01"),
    @("1982-83", "Moses Malone\malonmo01"),
    @("1981-82", "Moses Malone\malonmo01"),
    @("1980-81", "Julius Erving\ervinju01"),
    @("1979-80", "Kareem Abdul-Jabbar\abdulka01"),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
}

# --- Drop the Award column (col O); every row was "MVP" so it carried no
#     information once Season/Player became unique identifiers. ---
$ws.Columns("O").Delete() | Out-Null

# --- Restore default selection to A1-ish state used by the refreshed sheet ---
$ws.Range("O2").Select() | Out-Null

